$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.910.44'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +6.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.015.31'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +4.48%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.68'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.32'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +11.95%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.014.16'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +4.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.517'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.01'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.154'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.96%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.456'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +6.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000249'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +7.76%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.76'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +9.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.127'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.016.96'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +7.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.518.59'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +4.62%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.00'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +7.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.016.64'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '457.91'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +6.69%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +6.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.687'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +5.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.36'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +7.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.30'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.21%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +11.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.42'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.62'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +5.35%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.23'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +17.96%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +17.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0000104'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.92%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.37%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +5.47%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.44%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.996'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +4.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.78'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +7.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.14'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +12.13%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.05'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.52%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.99'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +7.16%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.122'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +7.28%  '
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.303'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +13.00%  '
$ws.Range('B43').Value = 'Arweave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '43.72'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +10.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.50'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.45%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.791.72'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.97%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '381.38'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +10.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0355'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +5.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '134.53'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.22%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.79'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +10.93%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +4.13%  '
